{"js": "// Replace each two-digit-division problem \"A\u00f7B=\" in the body (including\n// the practice table) with its corresponding new problem, in document\n// order. The document contains exactly 25 such problems (5 rows x 5\n// columns); map position i (0-based, document order) -> newValues[i].\nconst oldValues = [\n  \"15\u00f79=\", \"34\u00f77=\", \"41\u00f76=\", \"65\u00f76=\", \"17\u00f78=\",\n  \"16\u00f72=\", \"52\u00f76=\", \"52\u00f76=\", \"63\u00f74=\", \"64\u00f78=\",\n  \"44\u00f75=\", \"12\u00f72=\", \"59\u00f77=\", \"94\u00f76=\", \"70\u00f76=\",\n  \"23\u00f73=\", \"46\u00f75=\", \"80\u00f76=\", \"29\u00f75=\", \"46\u00f76=\",\n  \"44\u00f73=\", \"33\u00f73=\", \"64\u00f77=\", \"96\u00f76=\", \"58\u00f72=\"\n];\nconst newValues = [\n  \"12\u00f76=\", \"18\u00f75=\", \"99\u00f73=\", \"13\u00f77=\", \"68\u00f78=\",\n  \"45\u00f74=\", \"92\u00f76=\", \"75\u00f76=\", \"90\u00f72=\", \"96\u00f78=\",\n  \"61\u00f73=\", \"35\u00f74=\", \"59\u00f78=\", \"43\u00f74=\", \"48\u00f73=\",\n  \"30\u00f79=\", \"59\u00f76=\", \"39\u00f72=\", \"63\u00f76=\", \"66\u00f77=\",\n  \"31\u00f73=\", \"66\u00f79=\", \"38\u00f76=\", \"31\u00f79=\", \"43\u00f76=\"\n];\n\nconst body = context.document.body;\nconst results = body.search(\"\u00f7\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length !== oldValues.length) {\n  throw new Error(\n    `Expected ${oldValues.length} \"\u00f7\" matches, found ${results.items.length}`\n  );\n}\n\n// Expand each \"\u00f7\" hit to its enclosing paragraph so we replace the whole\n// \"A\u00f7B=\" run rather than just the division sign.\nconst paragraphs = results.items.map((r) => r.paragraphs.getFirst());\nfor (const p of paragraphs) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.length; i++) {\n  const actual = paragraphs[i].text;\n  if (actual !== oldValues[i]) {\n    throw new Error(\n      `Mismatch at problem #${i + 1}: expected \"${oldValues[i]}\", found \"${actual}\"`\n    );\n  }\n  paragraphs[i].getRange().insertText(newValues[i], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-division problem \"A\u00f7B=\" in the practice table\n# with its corresponding new problem, in document order. The document\n# contains exactly 25 such problems (5 rows x 5 columns); position i\n# (0-based, document order) maps oldValues[i] -> newValues[i]. Some old\n# values repeat (e.g. \"52\u00f76=\" appears twice), so matches must be walked\n# strictly left-to-right/top-to-bottom rather than replaced by a single\n# global Replace All.\n$oldValues = @(\n  \"15\u00f79=\", \"34\u00f77=\", \"41\u00f76=\", \"65\u00f76=\", \"17\u00f78=\",\n  \"16\u00f72=\", \"52\u00f76=\", \"52\u00f76=\", \"63\u00f74=\", \"64\u00f78=\",\n  \"44\u00f75=\", \"12\u00f72=\", \"59\u00f77=\", \"94\u00f76=\", \"70\u00f76=\",\n  \"23\u00f73=\", \"46\u00f75=\", \"80\u00f76=\", \"29\u00f75=\", \"46\u00f76=\",\n  \"44\u00f73=\", \"33\u00f73=\", \"64\u00f77=\", \"96\u00f76=\", \"58\u00f72=\"\n)\n$newValues = @(\n  \"12\u00f76=\", \"18\u00f75=\", \"99\u00f73=\", \"13\u00f77=\", \"68\u00f78=\",\n  \"45\u00f74=\", \"92\u00f76=\", \"75\u00f76=\", \"90\u00f72=\", \"96\u00f78=\",\n  \"61\u00f73=\", \"35\u00f74=\", \"59\u00f78=\", \"43\u00f74=\", \"48\u00f73=\",\n  \"30\u00f79=\", \"59\u00f76=\", \"39\u00f72=\", \"63\u00f76=\", \"66\u00f77=\",\n  \"31\u00f73=\", \"66\u00f79=\", \"38\u00f76=\", \"31\u00f79=\", \"43\u00f76=\"\n)\n\n$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Collapse(1)\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n  $found = $rng.Find.Execute($oldValues[$i], $false, $false, $false, $false, $false, $true, 1, $false, $newValues[$i], 1)\n  if (-not $found) {\n    throw \"Could not find occurrence #$($i + 1) (`\"$($oldValues[$i])`\") in document order\"\n  }\n  $rng.Collapse(0)\n}\n"}
